$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (rows 2-5) and old headers first
$ws.Range("A1:D5").ClearContents()

# New header row
$ws.Range("A1").Value = "project"
$ws.Range("B1").Value = "region"
$ws.Range("C1").Value = "zone"
$ws.Range("D1").Value = "name"

# New data row
$ws.Range("A2").Value = "mytemporaryproject28490"
$ws.Range("B2").Value = "us-central1"
$ws.Range("C2").Value = "us-central1-a"
$ws.Range("D2").Value = "sample-instance-200"

# Autofit columns based on content (matches bestFit widths in target XML)
$ws.Columns("A:D").AutoFit()

# Move the active selection
$ws.Range("D8").Select() | Out-Null
